$d = $word.ActiveDocument
$olds = @(
    "[[PERSON_34]] – „s [[PERSON_35]]“",
    "[[PERSON_36]] – „o [[PERSON_37]]“",
    "[[PERSON_38]] – „pro [[PERSON_38]]“",
    "[[PERSON_39]] – „s [[PERSON_39]]“",
    "[[PERSON_40]] – „k [[PERSON_41]]“",
    "[[PERSON_42]] – „s [[PERSON_42]]“",
    "[[PERSON_43]] – „o [[PERSON_43]]“",
    "[[PERSON_44]] – „u [[PERSON_45]]“",
    "[[PERSON_46]] – „k [[PERSON_46]]“",
    "[[PERSON_47]] – „se [[PERSON_47]]“",
    "[[PERSON_48]] – „u [[PERSON_48]]“",
    "[[PERSON_49]] – „o [[PERSON_50]]“",
    "[[PERSON_51]] – „s [[PERSON_51]]“",
    "[[PERSON_52]] – „k [[PERSON_53]]“",
    "[[PERSON_54]] – „od [[PERSON_55]]“",
    "[[PERSON_56]] – „s [[PERSON_56]]“",
    "[[PERSON_57]] – „u [[PERSON_58]]“",
    "[[PERSON_59]] – „o [[PERSON_60]]“",
    "[[PERSON_61]] – „k [[PERSON_62]]“",
    "V těchto řízeních bylo jednáno např. s [[PERSON_3]], [[PERSON_8]], [[PERSON_36]] či [[PERSON_63]].",
    "svědek [[PERSON_42]] (ve výpovědi označen jako „svědek Černého“),",
    "právní zástupkyně JUDr. [[PERSON_46]], advokátka,",
    "tlumočník [[PERSON_44]], zapsaný v seznamu tlumočníků.",
    "Neurologické testy č. NEU/2025/44119 provedené MUDr. [[PERSON_49]],",
    "Oční vyšetření č. OFT/2023/11281 provedené MUDr. [[PERSON_43]].",
    "Zvláštní pozornost byla věnována výsledkům [[PERSON_22]], [[PERSON_28]] a [[PERSON_59]].",
    "mobil [[PERSON_64]] S22, [[IMEI_1]],",
    "[[PERSON_65]] poskytly technické přístupy pro řešení kauz:",
    "právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_61]]),",
    "[[PERSON_56]] („výslech [[PERSON_56]]“),",
    "[[PERSON_51]] („výpověď [[PERSON_51]]“),",
    "[[PERSON_48]] („záznam o výslechu [[PERSON_48]]“),",
    "[[PERSON_27]] („výslech [[PERSON_66]]“).",
    "PhDr. [[PERSON_43]] – psychologický posudek,",
    "MUDr. [[PERSON_36]] – posudek z traumatologie,",
    "Tyto účty byly doloženy např. od [[PERSON_30]], [[PERSON_52]] nebo [[PERSON_67]].",
    "[[PERSON_54]],",
    "[[PERSON_63]],"
)
$news = @(
    "[[PERSON_34]] – „s [[PERSON_34]]“",
    "[[PERSON_35]] – „o [[PERSON_36]]“",
    "[[PERSON_37]] – „pro [[PERSON_37]]“",
    "[[PERSON_38]] – „s [[PERSON_38]]“",
    "[[PERSON_39]] – „k [[PERSON_40]]“",
    "[[PERSON_41]] – „s [[PERSON_41]]“",
    "[[PERSON_42]] – „o [[PERSON_42]]“",
    "[[PERSON_43]] – „u [[PERSON_44]]“",
    "[[PERSON_45]] – „k [[PERSON_45]]“",
    "[[PERSON_46]] – „se [[PERSON_46]]“",
    "[[PERSON_47]] – „u [[PERSON_47]]“",
    "[[PERSON_48]] – „o [[PERSON_49]]“",
    "[[PERSON_50]] – „s [[PERSON_50]]“",
    "[[PERSON_51]] – „k [[PERSON_52]]“",
    "[[PERSON_53]] – „od [[PERSON_54]]“",
    "[[PERSON_55]] – „s [[PERSON_55]]“",
    "[[PERSON_56]] – „u [[PERSON_57]]“",
    "[[PERSON_58]] – „o [[PERSON_59]]“",
    "[[PERSON_60]] – „k [[PERSON_61]]“",
    "V těchto řízeních bylo jednáno např. s [[PERSON_3]], [[PERSON_8]], [[PERSON_35]] či [[PERSON_62]].",
    "svědek [[PERSON_41]] (ve výpovědi označen jako „svědek Černého“),",
    "právní zástupkyně JUDr. [[PERSON_45]], advokátka,",
    "tlumočník [[PERSON_43]], zapsaný v seznamu tlumočníků.",
    "Neurologické testy č. NEU/2025/44119 provedené MUDr. [[PERSON_48]],",
    "Oční vyšetření č. OFT/2023/11281 provedené MUDr. [[PERSON_42]].",
    "Zvláštní pozornost byla věnována výsledkům [[PERSON_22]], [[PERSON_28]] a [[PERSON_58]].",
    "mobil [[PERSON_63]] S22, [[IMEI_1]],",
    "[[PERSON_64]] poskytly technické přístupy pro řešení kauz:",
    "právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_60]]),",
    "[[PERSON_55]] („výslech [[PERSON_55]]“),",
    "[[PERSON_50]] („výpověď [[PERSON_50]]“),",
    "[[PERSON_47]] („záznam o výslechu [[PERSON_47]]“),",
    "[[PERSON_27]] („výslech [[PERSON_65]]“).",
    "PhDr. [[PERSON_42]] – psychologický posudek,",
    "MUDr. [[PERSON_35]] – posudek z traumatologie,",
    "Tyto účty byly doloženy např. od [[PERSON_30]], [[PERSON_51]] nebo [[PERSON_66]].",
    "[[PERSON_53]],",
    "[[PERSON_62]],"
)

$notFound = @()
for ($i = 0; $i -lt $olds.Length; $i++) {
    $old = $olds[$i]
    $new = $news[$i]
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        $notFound = $notFound + $old
    }
}
if ($notFound.Length -gt 0) {
    Write-Output ("NOT FOUND COUNT=" + $notFound.Length)
    foreach ($nf in $notFound) {
        Write-Output ("MISSING: " + $nf)
    }
} else {
    Write-Output "ALL REPLACEMENTS APPLIED"
}
